$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Approved -> Rejected, and add a ReasonToReject of "Nil"
$ws.Range("I2").Value = "Rejected"
$ws.Range("J2").Value = "Nil"

# Row 3: Rejected -> Approved, clearing the now-inapplicable ReasonToReject
$ws.Range("I3").Value = "Approved"
$ws.Range("J3").ClearContents()

# Move the active selection to H16
$ws.Range("H16").Select()
